$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 04:35"

# --- Swap Belice / Nueva Caledonia (rows 193/194) including their D/H stats ---
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("A194").Value = "Belice"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0
$ws.Range("D194").Value = 16
$ws.Range("H194").Value = 2

# --- Swap Curazao / Dominica (rows 198/199) including their D/H stats ---
$ws.Range("A198").Value = "Dominica"
$ws.Range("A199").Value = "Curazao"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# --- Swap Sahara Occidental / San Bartolome (rows 215/216) ---
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Sahara Occidental"

# --- Updated country stats (Brasil, row 10) ---
$ws.Range("B10").Value = 178214
$ws.Range("C10").Value = 612
$ws.Range("E10").Value = 93156
$ws.Range("G10").Value = 57
$ws.Range("H10").Value = 12461

# --- Updated country stats (Chequia, row 51) ---
$ws.Range("B51").Value = 8221
$ws.Range("C51").Value = 23
$ws.Range("D51").Value = 4889
$ws.Range("E51").Value = 3049

# --- Updated country stats (Honduras, row 79) ---
$ws.Range("B79").Value = 2080
$ws.Range("E79").Value = 1758
